$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 281 (1-based), shifting existing rows 281..344 down to 282..345
$ws.Rows.Item(281).Insert()

# Populate the new row 281 with the new data point
$ws.Cells.Item(281, 1).Value = 11
$ws.Cells.Item(281, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(281, 3).Value = "Bíobío"
$ws.Cells.Item(281, 4).Value = 45244
$ws.Cells.Item(281, 5).Value = 8
$ws.Cells.Item(281, 6).Value = 100112003
$ws.Cells.Item(281, 7).Value = "Ajo"
$ws.Cells.Item(281, 8).Value = "Chino"
$ws.Cells.Item(281, 9).Value = "Primera"
$ws.Cells.Item(281, 10).Value = 250
$ws.Cells.Item(281, 11).Value = 22000
$ws.Cells.Item(281, 12).Value = 23000
$ws.Cells.Item(281, 13).Value = 22400
$ws.Cells.Item(281, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(281, 15).Value = "China"
$ws.Cells.Item(281, 16).Value = 2240
$ws.Cells.Item(281, 17).Value = 10
$ws.Cells.Item(281, 18).Value = "Hortaliza"
